$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SIQ")
$ws2 = $wb.Worksheets.Item("Version history")

# ---------------------------------------------------------------------------
# 1) Version history sheet: add new row 5 (version 1.3 entry)
# ---------------------------------------------------------------------------
# Copy style of the previous entry row (row4) into row5's date cell only, so
# the new date cell D5 picks up the same date number format (row A/B/C keep
# their existing style already).
$ws2.Range("D4").Copy($ws2.Range("D5"))

# New shared string #1: version-history comment (becomes shared string idx 63)
$ws2.Range("C5").Value = "Added questions about the login module and answers assumptions"

# ---------------------------------------------------------------------------
# 2) SIQ sheet: new login-module rows 19-22
# ---------------------------------------------------------------------------
# Copy formatting from the last existing data row (row18) down across the
# four new rows so the styles match (borders/fonts/number formats).
$ws1.Range("A18:E18").Copy($ws1.Range("A19:E22"))

# New IDs (shared strings idx 64-67)
$ws1.Range("A19").Value = "LH-SIQ-018"
$ws1.Range("A20").Value = "LH-SIQ-019"
$ws1.Range("A21").Value = "LH-SIQ-020"
$ws1.Range("A22").Value = "LH-SIQ-021"

# New feature label (shared string idx 68)
$ws1.Range("B19").Value = "Login"

# New question texts, entered in the same order the original author typed
# them (rows 19, 20, 22, then 21) so shared-string indices line up.
$ws1.Range("C19").Value = "Should the user login using the username and password?"
$ws1.Range("C20").Value = "Should only registered users be able to login?"
$ws1.Range("C22").Value = "Should we store registrants passwords using Hashing and SALT?"
$ws1.Range("C21").Value = "Should there be a generic error message if the user enters wrong email or password or leaves any of them empty"

# Remaining feature labels / answers / comments (reuse existing strings)
$ws1.Range("B20").Value = "Login"
$ws1.Range("B21").Value = "Login"
$ws1.Range("B22").Value = "Registration"
$ws1.Range("D19:D22").Value = "Yes"
$ws1.Range("E19:E22").Value = "no comment"

# Ensure row heights match the sheet's standard custom height.
$ws1.Range("A19:E22").RowHeight = 40.049999999999997

# ---------------------------------------------------------------------------
# 3) Resize Table1 to include the new rows
# ---------------------------------------------------------------------------
$lo = $ws1.ListObjects.Item(1)
$lo.Resize($ws1.Range("A1:E22"))

# ---------------------------------------------------------------------------
# 4) Extend data validations to the new rows
# ---------------------------------------------------------------------------
$dv1 = $ws1.Range("D2:D22")
$dv1.Validation.Delete()
$dv1.Validation.Add(3, 1, 1, '"Yes,No"')
$dv1.Validation.ErrorTitle = "Yes or NO"
$dv1.Validation.ShowInput = $false
$dv1.Validation.ShowError = $true

$dv2 = $ws1.Range("B2:C22")
$dv2.Validation.Delete()
$dv2.Validation.Add(0, 1, 0)
$dv2.Validation.InCellDropdown = $false
$dv2.Validation.ShowInput = $false
$dv2.Validation.ShowError = $false

# ---------------------------------------------------------------------------
# 5) Version history row5 remaining cells (version number + author + date)
# ---------------------------------------------------------------------------
$ws2.Range("A5").Value = 1.3
$ws2.Range("B5").Value = "Mahmoud Radi"
$ws2.Range("D5").Value = 45760

# ---------------------------------------------------------------------------
# 6) View state: selections + active sheet
# ---------------------------------------------------------------------------
# Set the selection on "Version history" first (it is the currently active
# sheet), then switch activation over to "SIQ" which becomes the active tab.
$ws2.Range("C9").Select()

$ws1.Activate()
$ws1.Range("C24").Select()
